$d = $word.ActiveDocument

# Font name Calibri -> Aptos, and sizes updated for the CLIN* custom styles
# (sz values below are in points = half-points/2)

$s = $d.Styles("CLIN1HEADING")
$s.Font.Name = "Aptos"
$s.Font.Size = 15

$s = $d.Styles("CLIN2SUBHEADINGS")
$s.Font.Name = "Aptos"
$s.Font.Size = 10

$s = $d.Styles("CLIN1HEADINGChar")
$s.Font.Name = "Aptos"
$s.Font.Size = 15

$s = $d.Styles("CLIN3BULLETPOINTS")
$s.Font.Name = "Aptos"
$s.Font.Size = 8

$s = $d.Styles("CLIN2SUBHEADINGSChar")
$s.Font.Name = "Aptos"
$s.Font.Size = 10

$s = $d.Styles("CLIN4")
$s.Font.Name = "Aptos"
$s.Font.Size = 5.5

$s = $d.Styles("CLIN3BULLETPOINTSChar")
$s.Font.Name = "Aptos"
$s.Font.Size = 8

$s = $d.Styles("CLIN4Char")
$s.Font.Name = "Aptos"
$s.Font.Size = 5.5
